$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The oldest reporting period (column D, "6 ماهه منتهی به 1399/06" / dated
# 1400-09-04 (3)) has rolled off the report; delete it entirely so every
# later column shifts one slot to the left (matches the target layout,
# incl. column widths and the row-21 string/number cut-over point).
$ws.Range("D1").EntireColumn.Delete()

# Bring up a brand-new rightmost column (M) for the newest period by
# cloning column L's formatting/styles, then overwrite with the new data.
$ws.Range("L1:L28").Copy($ws.Range("M1"))
$ws.Range("M1").ColumnWidth = 28.1

# The string that used to read "1401-11-01 (6)" slid from J9 to I9 during
# the shift; the newest publish-date revision renames it.
$ws.Range("I9").Value = "1402-02-27 (7)"

# New column M: header + publish date for the 12-month period ending 1401/12.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-27"

# New column M figures (dollar cumulative income statement).
$ws.Range("M11").Value = 57825
$ws.Range("M12").Value = -46116
$ws.Range("M13").Value = 11709
$ws.Range("M14").Value = -2791
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 57
$ws.Range("M17").Value = 8975
$ws.Range("M18").Value = -1234
$ws.Range("M19").Value = -157
$ws.Range("M20").Value = 7584
$ws.Range("M21").Value = "-"
$ws.Range("M22").Value = 7584
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 7584
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 1984
$ws.Range("M27").Value = 0
